$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vars")

# --- Header row (row 1): new columns J="mult", K="show", L="rtol", M="atol" ---

# J1, K1 use the same format as the existing "tex"/"html" header cells (H1/I1)
$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)
$ws.Range("J1").Value = "mult"
$ws.Range("K1").Value = "show"

# L1, M1 use a bold header style too, but with an explicit black font color
# (this is a new, slightly different header style introduced for the two
# new tolerance columns)
$ws.Range("I1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)
$ws.Range("L1:M1").Font.Color = 0
$ws.Range("L1").Value = "rtol"
$ws.Range("M1").Value = "atol"

# --- Data rows 2-5 ---

# "mult" and "show" are plain integers (0/1 flags), formatted like the other
# plain numeric columns (D5:F5 use the default numeric style)
$ws.Range("D5").Copy()
$ws.Range("J2:K5").PasteSpecial(-4122)
$ws.Range("J2:J5").Value = 1
$ws.Range("K2").Value = 0
$ws.Range("K3:K5").Value = 1

# "rtol" and "atol" are small tolerances, formatted in scientific notation
# like the other scientific-format columns (D3:F3 use numFmt 0.00E+000)
$ws.Range("D3").Copy()
$ws.Range("L2:M5").PasteSpecial(-4122)
$ws.Range("L2:M5").Value = 0.000001

# --- Make "vars" the active sheet / active tab, with M3 selected ---
# (previously "funs" was the active tab with G21 selected on "vars")
$ws.Activate()
$ws.Range("M3").Select()
